$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = -19.02791139457464
$ws.Cells.Item(2, 3).Value = -0.7396535990740564
$ws.Cells.Item(2, 4).Value = -19.02791139457464
$ws.Cells.Item(2, 5).Value = -19.02791139457464
$ws.Cells.Item(2, 6).Value = -19.02791139457464
$ws.Cells.Item(2, 7).Value = -19.02791139457464
$ws.Cells.Item(2, 8).Value = -19.02791139457464
$ws.Cells.Item(2, 9).Value = -19.02791139457464
$ws.Cells.Item(2, 10).Value = -19.02791139457464
$ws.Cells.Item(2, 11).Value = -19.02791139457464

$ws.Cells.Item(3, 2).Value = -19.02791139457464
$ws.Cells.Item(3, 3).Value = -19.02791139457464
$ws.Cells.Item(3, 4).Value = -19.02791139457464
$ws.Cells.Item(3, 5).Value = -19.02791139457464
$ws.Cells.Item(3, 6).Value = -19.02791139457464
$ws.Cells.Item(3, 7).Value = -19.02791139457464
$ws.Cells.Item(3, 8).Value = -19.02791139457464
$ws.Cells.Item(3, 9).Value = -0.1706441260524985
$ws.Cells.Item(3, 10).Value = -19.02791139457464
$ws.Cells.Item(3, 11).Value = -19.02791139457464

$ws.Cells.Item(4, 2).Value = -19.02791139457464
$ws.Cells.Item(4, 3).Value = -0.5673075182622004
$ws.Cells.Item(4, 4).Value = 0.08614862496725685
$ws.Cells.Item(4, 5).Value = -19.02791139457464
$ws.Cells.Item(4, 6).Value = 4.009362235884019
$ws.Cells.Item(4, 7).Value = -19.02791139457464
$ws.Cells.Item(4, 8).Value = 2.356864452270391
$ws.Cells.Item(4, 9).Value = -19.02791139457464
$ws.Cells.Item(4, 10).Value = 3.101371936852857
$ws.Cells.Item(4, 11).Value = -19.02791139457464

$ws.Cells.Item(5, 2).Value = -19.02791139457464
$ws.Cells.Item(5, 3).Value = -0.0413665773124736
$ws.Cells.Item(5, 4).Value = -19.02791139457464
$ws.Cells.Item(5, 5).Value = -19.02791139457464
$ws.Cells.Item(5, 6).Value = -19.02791139457464
$ws.Cells.Item(5, 7).Value = 4.321925530837624
$ws.Cells.Item(5, 8).Value = -19.02791139457464
$ws.Cells.Item(5, 9).Value = -19.02791139457464
$ws.Cells.Item(5, 10).Value = -19.02791139457464
$ws.Cells.Item(5, 11).Value = -19.02791139457464

$ws.Cells.Item(6, 2).Value = -19.02791139457464
$ws.Cells.Item(6, 3).Value = -19.02791139457464
$ws.Cells.Item(6, 4).Value = -19.02791139457464
$ws.Cells.Item(6, 5).Value = -19.02791139457464
$ws.Cells.Item(6, 6).Value = -19.02791139457464
$ws.Cells.Item(6, 7).Value = -19.02791139457464
$ws.Cells.Item(6, 8).Value = -19.02791139457464
$ws.Cells.Item(6, 9).Value = -19.02791139457464
$ws.Cells.Item(6, 10).Value = -19.02791139457464
$ws.Cells.Item(6, 11).Value = -19.02791139457464

$ws.Cells.Item(7, 2).Value = 3.24018080132918
$ws.Cells.Item(7, 3).Value = -19.02791139457464
$ws.Cells.Item(7, 4).Value = -19.02791139457464
$ws.Cells.Item(7, 5).Value = -19.02791139457464
$ws.Cells.Item(7, 6).Value = -19.02791139457464
$ws.Cells.Item(7, 7).Value = -19.02791139457464
$ws.Cells.Item(7, 8).Value = -19.02791139457464
$ws.Cells.Item(7, 9).Value = -19.02791139457464
$ws.Cells.Item(7, 10).Value = -19.02791139457464
$ws.Cells.Item(7, 11).Value = -19.02791139457464

$ws.Cells.Item(8, 2).Value = -19.02791139457464
$ws.Cells.Item(8, 3).Value = -19.02791139457464
$ws.Cells.Item(8, 4).Value = -19.02791139457464
$ws.Cells.Item(8, 5).Value = 1.519112081480545
$ws.Cells.Item(8, 6).Value = -19.02791139457464
$ws.Cells.Item(8, 7).Value = -19.02791139457464
$ws.Cells.Item(8, 8).Value = -19.02791139457464
$ws.Cells.Item(8, 9).Value = -19.02791139457464
$ws.Cells.Item(8, 10).Value = -19.02791139457464
$ws.Cells.Item(8, 11).Value = -19.02791139457464

$ws.Cells.Item(9, 2).Value = 3.399286084431158
$ws.Cells.Item(9, 3).Value = -19.02791139457464
$ws.Cells.Item(9, 4).Value = -19.02791139457464
$ws.Cells.Item(9, 5).Value = -19.02791139457464
$ws.Cells.Item(9, 6).Value = -19.02791139457464
$ws.Cells.Item(9, 7).Value = -19.02791139457464
$ws.Cells.Item(9, 8).Value = -19.02791139457464
$ws.Cells.Item(9, 9).Value = -19.02791139457464
$ws.Cells.Item(9, 10).Value = -19.02791139457464
$ws.Cells.Item(9, 11).Value = -19.02791139457464

$ws.Cells.Item(10, 2).Value = -19.02791139457464
$ws.Cells.Item(10, 3).Value = -19.02791139457464
$ws.Cells.Item(10, 4).Value = -19.02791139457464
$ws.Cells.Item(10, 5).Value = -19.02791139457464
$ws.Cells.Item(10, 6).Value = -19.02791139457464
$ws.Cells.Item(10, 7).Value = -19.02791139457464
$ws.Cells.Item(10, 8).Value = -19.02791139457464
$ws.Cells.Item(10, 9).Value = 0.2447571142019817
$ws.Cells.Item(10, 10).Value = -19.02791139457464
$ws.Cells.Item(10, 11).Value = 1.90450883474331

$ws.Cells.Item(11, 2).Value = -19.02791139457464
$ws.Cells.Item(11, 3).Value = -19.02791139457464
$ws.Cells.Item(11, 4).Value = -19.02791139457464
$ws.Cells.Item(11, 5).Value = 2.389465750346117
$ws.Cells.Item(11, 6).Value = -19.02791139457464
$ws.Cells.Item(11, 7).Value = -19.02791139457464
$ws.Cells.Item(11, 8).Value = -19.02791139457464
$ws.Cells.Item(11, 9).Value = -19.02791139457464
$ws.Cells.Item(11, 10).Value = -19.02791139457464
$ws.Cells.Item(11, 11).Value = 1.294369978212916

$ws.Cells.Item(12, 2).Value = -19.02791139457464
$ws.Cells.Item(12, 3).Value = -19.02791139457464
$ws.Cells.Item(12, 4).Value = -19.02791139457464
$ws.Cells.Item(12, 5).Value = -19.02791139457464
$ws.Cells.Item(12, 6).Value = -19.02791139457464
$ws.Cells.Item(12, 7).Value = -19.02791139457464
$ws.Cells.Item(12, 8).Value = -19.02791139457464
$ws.Cells.Item(12, 9).Value = -19.02791139457464
$ws.Cells.Item(12, 10).Value = -19.02791139457464
$ws.Cells.Item(12, 11).Value = -19.02791139457464

$ws.Cells.Item(13, 2).Value = -19.02791139457464
$ws.Cells.Item(13, 3).Value = -19.02791139457464
$ws.Cells.Item(13, 4).Value = -19.02791139457464
$ws.Cells.Item(13, 5).Value = 1.920924199939891
$ws.Cells.Item(13, 6).Value = -19.02791139457464
$ws.Cells.Item(13, 7).Value = -19.02791139457464
$ws.Cells.Item(13, 8).Value = -19.02791139457464
$ws.Cells.Item(13, 9).Value = -19.02791139457464
$ws.Cells.Item(13, 10).Value = 0.5398683277989161
$ws.Cells.Item(13, 11).Value = 2.992328318578271

$ws.Cells.Item(14, 2).Value = -19.02791139457464
$ws.Cells.Item(14, 3).Value = -19.02791139457464
$ws.Cells.Item(14, 4).Value = 1.16902325457858
$ws.Cells.Item(14, 5).Value = -19.02791139457464
$ws.Cells.Item(14, 6).Value = -19.02791139457464
$ws.Cells.Item(14, 7).Value = -19.02791139457464
$ws.Cells.Item(14, 8).Value = -19.02791139457464
$ws.Cells.Item(14, 9).Value = -19.02791139457464
$ws.Cells.Item(14, 10).Value = -19.02791139457464
$ws.Cells.Item(14, 11).Value = 1.290601626401876

$ws.Cells.Item(15, 2).Value = -19.02791139457464
$ws.Cells.Item(15, 3).Value = -19.02791139457464
$ws.Cells.Item(15, 4).Value = -0.5162320817279584
$ws.Cells.Item(15, 5).Value = -19.02791139457464
$ws.Cells.Item(15, 6).Value = -19.02791139457464
$ws.Cells.Item(15, 7).Value = -19.02791139457464
$ws.Cells.Item(15, 8).Value = -19.02791139457464
$ws.Cells.Item(15, 9).Value = -19.02791139457464
$ws.Cells.Item(15, 10).Value = -19.02791139457464
$ws.Cells.Item(15, 11).Value = -19.02791139457464

$ws.Cells.Item(16, 2).Value = -19.02791139457464
$ws.Cells.Item(16, 3).Value = -19.02791139457464
$ws.Cells.Item(16, 4).Value = -19.02791139457464
$ws.Cells.Item(16, 5).Value = -19.02791139457464
$ws.Cells.Item(16, 6).Value = -19.02791139457464
$ws.Cells.Item(16, 7).Value = -19.02791139457464
$ws.Cells.Item(16, 8).Value = -19.02791139457464
$ws.Cells.Item(16, 9).Value = -19.02791139457464
$ws.Cells.Item(16, 10).Value = 2.059329349395679
$ws.Cells.Item(16, 11).Value = -19.02791139457464

$ws.Cells.Item(17, 2).Value = -19.02791139457464
$ws.Cells.Item(17, 3).Value = 0.4262774336570801
$ws.Cells.Item(17, 4).Value = -0.4563476017375718
$ws.Cells.Item(17, 5).Value = -19.02791139457464
$ws.Cells.Item(17, 6).Value = -19.02791139457464
$ws.Cells.Item(17, 7).Value = -19.02791139457464
$ws.Cells.Item(17, 8).Value = 2.21827308122441
$ws.Cells.Item(17, 9).Value = 0.2112978625204684
$ws.Cells.Item(17, 10).Value = 1.754809329085983
$ws.Cells.Item(17, 11).Value = -19.02791139457464

$ws.Cells.Item(18, 2).Value = -19.02791139457464
$ws.Cells.Item(18, 3).Value = -19.02791139457464
$ws.Cells.Item(18, 4).Value = -19.02791139457464
$ws.Cells.Item(18, 5).Value = -19.02791139457464
$ws.Cells.Item(18, 6).Value = -19.02791139457464
$ws.Cells.Item(18, 7).Value = -19.02791139457464
$ws.Cells.Item(18, 8).Value = 1.858464686038769
$ws.Cells.Item(18, 9).Value = -1.397921628160154
$ws.Cells.Item(18, 10).Value = 1.275629690023351
$ws.Cells.Item(18, 11).Value = -19.02791139457464

$ws.Cells.Item(19, 2).Value = -19.02791139457464
$ws.Cells.Item(19, 3).Value = -19.02791139457464
$ws.Cells.Item(19, 4).Value = 2.879074957795956
$ws.Cells.Item(19, 5).Value = -19.02791139457464
$ws.Cells.Item(19, 6).Value = -19.02791139457464
$ws.Cells.Item(19, 7).Value = -19.02791139457464
$ws.Cells.Item(19, 8).Value = 1.849565220866906
$ws.Cells.Item(19, 9).Value = 1.677502473001379
$ws.Cells.Item(19, 10).Value = -19.02791139457464
$ws.Cells.Item(19, 11).Value = -19.02791139457464

$ws.Cells.Item(20, 2).Value = -19.02791139457464
$ws.Cells.Item(20, 3).Value = 3.286794921096016
$ws.Cells.Item(20, 4).Value = 2.982772752077301
$ws.Cells.Item(20, 5).Value = -19.02791139457464
$ws.Cells.Item(20, 6).Value = 1.961918940414754
$ws.Cells.Item(20, 7).Value = -19.02791139457464
$ws.Cells.Item(20, 8).Value = 0.6008278181464697
$ws.Cells.Item(20, 9).Value = 3.721448429944212
$ws.Cells.Item(20, 10).Value = -19.02791139457464
$ws.Cells.Item(20, 11).Value = 1.765359300869269

$ws.Cells.Item(21, 2).Value = -19.02791139457464
$ws.Cells.Item(21, 3).Value = 2.733636997379394
$ws.Cells.Item(21, 4).Value = -19.02791139457464
$ws.Cells.Item(21, 5).Value = 3.019263452114498
$ws.Cells.Item(21, 6).Value = -19.02791139457464
$ws.Cells.Item(21, 7).Value = -19.02791139457464
$ws.Cells.Item(21, 8).Value = 0.5629990588976846
$ws.Cells.Item(21, 9).Value = -19.02791139457464
$ws.Cells.Item(21, 10).Value = -19.02791139457464
$ws.Cells.Item(21, 11).Value = -19.02791139457464

Write-Output "applied changes"